$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = "26.724.77"
$ws.Range("D3").Value = "1.600.96"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue $ws.Range("D5") "211.50"
$ws.Range("E5").Value = "  +0.05%  "
Set-TextValue $ws.Range("D6") "0.514"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +0.15%  "
Set-TextValue $ws.Range("D8") "0.0618"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.38%  "
Set-TextValue $ws.Range("D10") "19.69"
$ws.Range("E10").Value = "  +0.64%  "
Set-TextValue $ws.Range("D11") "0.0846"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.825.58"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.589.03"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  -0.02%  "
Set-TextValue $ws.Range("D16") "65.02"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "0.0₃0738"
$ws.Range("E17").Value = "  +0.06%  "
Set-TextValue $ws.Range("D18") "209.80"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +0.16%  "
Set-TextValue $ws.Range("D20") "7.16"
$ws.Range("E20").Value = "  +2.06%  "
Set-TextValue $ws.Range("D22") "2.24"
$ws.Range("E22").Value = "  -3.92%  "
Set-TextValue $ws.Range("D23") "8.98"
$ws.Range("E23").Value = "  -0.05%  "
Set-TextValue $ws.Range("D24") "143.45"
$ws.Range("E24").Value = "  -0.52%  "
Set-TextValue $ws.Range("D25") "1.01"
$ws.Range("E25").Value = "  +0.14%  "
Set-TextValue $ws.Range("D26") "7.07"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  -0.67%  "
Set-TextValue $ws.Range("D28") "15.32"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.289.00"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  +11.28%  "
$ws.Range("E38").Value = "  -0.07%  "
Set-TextValue $ws.Range("D39") "0.833"
$ws.Range("E39").Value = "  -0.09%  "
Set-TextValue $ws.Range("D40") "5.41"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  -0.60%  "
Set-TextValue $ws.Range("D42") "0.780"
$ws.Range("E42").Value = "  -0.20%  "
Set-TextValue $ws.Range("D43") "62.83"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").Value = "1.736.90"
Set-TextValue $ws.Range("D45") "90.48"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D48") "0.102"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0516"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws.Range("D50") "1.00"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.47"
$ws.Range("E51").Value = "  +1.07%  "
